# edit.ps1 - apply the "Week 7 -> Week 6" / new bullet / run-merge changes
# described by the OOXML diff for modules/week06/slides-06.pptx.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 ("title" slide): the sub-header runs "Week " + "7" collapse into a
# single run reading "Week 6".
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$headerShape = $slide1.Shapes.Item(4)
$headerText = $headerShape.TextFrame.TextRange
$weekPara = $headerText.Paragraphs(2, 1)

# Round-trip through a placeholder value first so the engine actually
# rewrites the paragraph's runs (assigning the exact same text the
# paragraph already renders as is treated as a no-op otherwise).
$weekPara.Text = "__TMP__"
$headerText.Paragraphs(2, 1).Text = "Week 6"

# ---------------------------------------------------------------------------
# Slide 3 ("Programming with databases... when?"): add a new sub-bullet and
# fix up a split run.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$bodyShape = $slide3.Shapes.Item(2)
$bodyText = $bodyShape.TextFrame.TextRange

# Insert "Create cool data visualizations" as a new lvl-1 bullet right after
# "Alternative to exporting data as CSV" (paragraph 2), before the blank
# lvl-1 paragraph that precedes "Support application".
$altPara = $bodyText.Paragraphs(2, 1)
$altPara.InsertAfter("`rCreate cool data visualizations")

# Re-resolve the body range/paragraph indices after the insertion above -
# "E.g" + "., Shiny" (two runs) now live in paragraph 6.
$bodyText = $bodyShape.TextFrame.TextRange
$egPara = $bodyText.Paragraphs(6, 1)
$egPara.Text = "__TMP__"
$bodyShape.TextFrame.TextRange.Paragraphs(6, 1).Text = "E.g., Shiny"
